$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '28.107.86'
$cell.Style = $origStyle
$cell = $ws.Range('E2')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.61%  '
$cell.Style = $origStyle
$cell = $ws.Range('D3')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.820.29'
$cell.Style = $origStyle
$cell = $ws.Range('E3')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.63%  '
$cell.Style = $origStyle
$cell = $ws.Range('E4')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.22%  '
$cell.Style = $origStyle
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '310.62'
$cell.Style = $origStyle
$cell = $ws.Range('E5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.04%  '
$cell.Style = $origStyle
$cell = $ws.Range('E6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.25%  '
$cell.Style = $origStyle
$cell = $ws.Range('D7')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.4977'
$cell.Style = $origStyle
$cell = $ws.Range('E7')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.55%  '
$cell.Style = $origStyle
$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.3913'
$cell.Style = $origStyle
$cell = $ws.Range('E8')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.75%  '
$cell.Style = $origStyle
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.09897'
$cell.Style = $origStyle
$cell = $ws.Range('E9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +25.61%  '
$cell.Style = $origStyle
$cell = $ws.Range('E10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.05%  '
$cell.Style = $origStyle
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '40.94'
$cell.Style = $origStyle
$cell = $ws.Range('E11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.40%  '
$cell.Style = $origStyle
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.450'
$cell.Style = $origStyle
$cell = $ws.Range('E12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +1.66%  '
$cell.Style = $origStyle
$cell = $ws.Range('D13')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '20.61'
$cell.Style = $origStyle
$cell = $ws.Range('E13')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.87%  '
$cell.Style = $origStyle
$cell = $ws.Range('E14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.27%  '
$cell.Style = $origStyle
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.814.44'
$cell.Style = $origStyle
$cell = $ws.Range('E15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +1.08%  '
$cell.Style = $origStyle
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.291'
$cell.Style = $origStyle
$cell = $ws.Range('E16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.52%  '
$cell.Style = $origStyle
$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.00001143'
$cell.Style = $origStyle
$cell = $ws.Range('E17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +5.67%  '
$cell.Style = $origStyle
$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '92.39'
$cell.Style = $origStyle
$cell = $ws.Range('E18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.44%  '
$cell.Style = $origStyle
$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.06646'
$cell.Style = $origStyle
$cell = $ws.Range('E19')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +1.14%  '
$cell.Style = $origStyle
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = $origStyle
$cell = $ws.Range('E20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.27%  '
$cell.Style = $origStyle
$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '17.21'
$cell.Style = $origStyle
$cell = $ws.Range('E21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.66%  '
$cell.Style = $origStyle
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.983'
$cell.Style = $origStyle
$cell = $ws.Range('E22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.45%  '
$cell.Style = $origStyle
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '28.157.80'
$cell.Style = $origStyle
$cell = $ws.Range('E23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.57%  '
$cell.Style = $origStyle
$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '11.27'
$cell.Style = $origStyle
$cell = $ws.Range('E24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +1.07%  '
$cell.Style = $origStyle
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.232'
$cell.Style = $origStyle
$cell = $ws.Range('E25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.31%  '
$cell.Style = $origStyle
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '159.23'
$cell.Style = $origStyle
$cell = $ws.Range('E26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.99%  '
$cell.Style = $origStyle
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '20.78'
$cell.Style = $origStyle
$cell = $ws.Range('E27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +1.00%  '
$cell.Style = $origStyle
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.022.83'
$cell.Style = $origStyle
$cell = $ws.Range('E28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.62%  '
$cell.Style = $origStyle
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.413'
$cell.Style = $origStyle
$cell = $ws.Range('E29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.33%  '
$cell.Style = $origStyle
$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '126.73'
$cell.Style = $origStyle
$cell = $ws.Range('E30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.85%  '
$cell.Style = $origStyle
$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.1054'
$cell.Style = $origStyle
$cell = $ws.Range('E31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.02%  '
$cell.Style = $origStyle
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.038'
$cell.Style = $origStyle
$cell = $ws.Range('E32')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.05%  '
$cell.Style = $origStyle
$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.579'
$cell.Style = $origStyle
$cell = $ws.Range('E33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.08%  '
$cell.Style = $origStyle
$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.606'
$cell.Style = $origStyle
$cell = $ws.Range('E34')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.63%  '
$cell.Style = $origStyle
$cell = $ws.Range('D35')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.06702'
$cell.Style = $origStyle
$cell = $ws.Range('E35')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -6.48%  '
$cell.Style = $origStyle
$cell = $ws.Range('E36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.60%  '
$cell.Style = $origStyle
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '8.912'
$cell.Style = $origStyle
$cell = $ws.Range('E37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.24%  '
$cell.Style = $origStyle
$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.2145'
$cell.Style = $origStyle
$cell = $ws.Range('E38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.50%  '
$cell.Style = $origStyle
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '11.42'
$cell.Style = $origStyle
$cell = $ws.Range('E39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.11%  '
$cell.Style = $origStyle
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.971'
$cell.Style = $origStyle
$cell = $ws.Range('E40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.76%  '
$cell.Style = $origStyle
$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.6214'
$cell.Style = $origStyle
$cell = $ws.Range('E41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.03%  '
$cell.Style = $origStyle
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.179'
$cell.Style = $origStyle
$cell = $ws.Range('E42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +2.15%  '
$cell.Style = $origStyle
$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = $origStyle
$cell = $ws.Range('E43')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.30%  '
$cell.Style = $origStyle
$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '13.15'
$cell.Style = $origStyle
$cell = $ws.Range('E44')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -0.60%  '
$cell.Style = $origStyle
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.5908'
$cell.Style = $origStyle
$cell = $ws.Range('E45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.14%  '
$cell.Style = $origStyle
$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.701'
$cell.Style = $origStyle
$cell = $ws.Range('E46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.24%  '
$cell.Style = $origStyle
$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.275'
$cell.Style = $origStyle
$cell = $ws.Range('E47')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -3.63%  '
$cell.Style = $origStyle
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '124.19'
$cell.Style = $origStyle
$cell = $ws.Range('E48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.10%  '
$cell.Style = $origStyle
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.946'
$cell.Style = $origStyle
$cell = $ws.Range('E49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  +0.23%  '
$cell.Style = $origStyle
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.181'
$cell.Style = $origStyle
$cell = $ws.Range('E50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -2.51%  '
$cell.Style = $origStyle
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.06770'
$cell.Style = $origStyle
$cell = $ws.Range('E51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '  -1.40%  '
$cell.Style = $origStyle
